# Apply edit: split the last paragraph (emoji paragraph), remove its paragraph-mark
# rFonts hint, and append new commentary paragraphs ("first part of hw3").
$d = $word.ActiveDocument

# Locate the target paragraph: the one that currently ends the document, containing
# a tab followed by "其他該用view的地方我都有用view" + the emoji AlternateContent run.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*view*view*") {
        $targetPara = $cand
    }
}

$insertXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex"><w:body><w:p><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>其他該用</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>view</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>的地方我都有用</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>view</w:t></w:r><w:r><w:rPr><mc:AlternateContent><mc:Choice Requires="w16se"/><mc:Fallback><w:rFonts w:ascii="Segoe UI Emoji" w:eastAsia="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/></mc:Fallback></mc:AlternateContent></w:rPr><mc:AlternateContent><mc:Choice Requires="w16se"><w16se:symEx w16se:font="Segoe UI Emoji" w16se:char="1F60A"/></mc:Choice><mc:Fallback><w:t>😊</w:t></mc:Fallback></mc:AlternateContent></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>備註</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>我使用</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>tablet</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>去弄介面的時候，他跑出來的不是</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>xlarge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>而是</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>sw600dp</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>然而課程內容照做就長那樣，希望不要太過刁難我，謝謝</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:ind w:firstLine="480"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>我也有用</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">nexus 10 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>去看過，所以就是那個</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>sw600dp</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>ZZDUMMYZZ</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$full = $d.Range($targetPara.Range.Start, $d.Content.End)
$full.InsertXML($insertXml)

# The inserted package ends with a throwaway paragraph (ZZDUMMYZZ) so that Word performs a
# clean full-range replace instead of merging into the original paragraph's trailing
# bookmark/run. Remove that helper paragraph mark now.
$lastIndex = $d.Paragraphs.Count
$prevEnd = $d.Paragraphs.Item($lastIndex - 1).Range.End
$cleanup = $d.Range($prevEnd - 1, $d.Content.End)
$cleanup.Delete()

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
